$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.226.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.559.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.43%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.559.35"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.28%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +4.01%  "
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.124.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.56%  "
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.559.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.347.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.09%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.608"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.703.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +8.68%  "
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.161"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.556.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "USDe"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("E38").Value = "  +4.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0859"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("E47").Value = "  +4.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.91%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.01%  "
